$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "-6.6559480162360725"
$ws.Range("C2").Style = "Normal"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "0.05868961271903217"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "53.01895184036078"
$ws.Range("E2").Style = "Normal"
$ws.Range("R2").NumberFormat = "@"
$ws.Range("R2").Value = "-75.28451095800045"
$ws.Range("R2").Style = "Normal"
$ws.Range("S2").NumberFormat = "@"
$ws.Range("S2").Value = "-41.154313722205565"
$ws.Range("S2").Style = "Normal"
$ws.Range("T2").NumberFormat = "@"
$ws.Range("T2").Value = "41.27169294764363"
$ws.Range("T2").Style = "Normal"
$ws.Range("U2").NumberFormat = "@"
$ws.Range("U2").Value = "48.54333966761809"
$ws.Range("U2").Style = "Normal"
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "-0.19837811831062346"
$ws.Range("C3").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.4192536029363474"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "42.68215859751152"
$ws.Range("E3").Style = "Normal"
$ws.Range("R3").NumberFormat = "@"
$ws.Range("R3").Value = "-46.27987070385327"
$ws.Range("R3").Style = "Normal"
$ws.Range("S3").NumberFormat = "@"
$ws.Range("S3").Value = "42.17787155491847"
$ws.Range("S3").Style = "Normal"
$ws.Range("T3").NumberFormat = "@"
$ws.Range("T3").Value = "42.64785102473808"
$ws.Range("T3").Style = "Normal"
$ws.Range("U3").NumberFormat = "@"
$ws.Range("U3").Value = "-39.33936434904577"
$ws.Range("U3").Style = "Normal"
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "33.8602244806344"
$ws.Range("C4").Style = "Normal"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "51.395636417917046"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "44.866183537550796"
$ws.Range("E4").Style = "Normal"
$ws.Range("R4").NumberFormat = "@"
$ws.Range("R4").Value = "74.22259317996169"
$ws.Range("R4").Style = "Normal"
$ws.Range("S4").NumberFormat = "@"
$ws.Range("S4").Value = "-41.57296809325818"
$ws.Range("S4").Style = "Normal"
$ws.Range("T4").NumberFormat = "@"
$ws.Range("T4").Value = "59.064697911019955"
$ws.Range("T4").Style = "Normal"
$ws.Range("U4").NumberFormat = "@"
$ws.Range("U4").Value = "43.72657492481414"
$ws.Range("U4").Style = "Normal"
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "21.993469305584206"
$ws.Range("C5").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "46.80857627733623"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "48.17201427980745"
$ws.Range("E5").Style = "Normal"
$ws.Range("R5").NumberFormat = "@"
$ws.Range("R5").Value = "55.55192158709969"
$ws.Range("R5").Style = "Normal"
$ws.Range("S5").NumberFormat = "@"
$ws.Range("S5").Value = "-61.19519691943532"
$ws.Range("S5").Style = "Normal"
$ws.Range("T5").NumberFormat = "@"
$ws.Range("T5").Value = "45.38224124859202"
$ws.Range("T5").Style = "Normal"
$ws.Range("U5").NumberFormat = "@"
$ws.Range("U5").Value = "48.234911306080434"
$ws.Range("U5").Style = "Normal"
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "-7.128249775739491"
$ws.Range("C6").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "-5.019466295626227"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "46.58977291779643"
$ws.Range("E6").Style = "Normal"
$ws.Range("R6").NumberFormat = "@"
$ws.Range("R6").Value = "-63.48833772415064"
$ws.Range("R6").Style = "Normal"
$ws.Range("S6").NumberFormat = "@"
$ws.Range("S6").Value = "45.01427121244513"
$ws.Range("S6").Style = "Normal"
$ws.Range("T6").NumberFormat = "@"
$ws.Range("T6").Value = "-42.290311590839494"
$ws.Range("T6").Style = "Normal"
$ws.Range("U6").NumberFormat = "@"
$ws.Range("U6").Value = "32.25137899958704"
$ws.Range("U6").Style = "Normal"
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "-35.72793970228871"
$ws.Range("C7").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "-64.3733722638362"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "50.78910541833564"
$ws.Range("E7").Style = "Normal"
$ws.Range("R7").NumberFormat = "@"
$ws.Range("R7").Value = "-78.44800770501213"
$ws.Range("R7").Style = "Normal"
$ws.Range("S7").NumberFormat = "@"
$ws.Range("S7").Value = "35.63756086198221"
$ws.Range("S7").Style = "Normal"
$ws.Range("T7").NumberFormat = "@"
$ws.Range("T7").Value = "-64.3733722638362"
$ws.Range("T7").Style = "Normal"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "6.099285711554951"
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "39.49336304668445"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "44.031188642542645"
$ws.Range("E8").Style = "Normal"
$ws.Range("R8").NumberFormat = "@"
$ws.Range("R8").Value = "-49.97706823054826"
$ws.Range("R8").Style = "Normal"
$ws.Range("S8").NumberFormat = "@"
$ws.Range("S8").Value = "-45.544626203098446"
$ws.Range("S8").Style = "Normal"
$ws.Range("T8").NumberFormat = "@"
$ws.Range("T8").Value = "39.49336304668445"
$ws.Range("T8").Style = "Normal"
$ws.Range("U8").NumberFormat = "@"
$ws.Range("U8").Value = "44.79880536049031"
$ws.Range("U8").Style = "Normal"
$ws.Range("V8").NumberFormat = "@"
$ws.Range("V8").Value = "41.72595458424669"
$ws.Range("V8").Style = "Normal"
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "-13.297306885160703"
$ws.Range("C9").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "-12.320497346273104"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "47.91974240997002"
$ws.Range("E9").Style = "Normal"
$ws.Range("R9").NumberFormat = "@"
$ws.Range("R9").Value = "30.727605220104557"
$ws.Range("R9").Style = "Normal"
$ws.Range("S9").NumberFormat = "@"
$ws.Range("S9").Value = "-71.04716843815358"
$ws.Range("S9").Style = "Normal"
$ws.Range("T9").NumberFormat = "@"
$ws.Range("T9").Value = "-52.65798943431227"
$ws.Range("T9").Style = "Normal"
$ws.Range("U9").NumberFormat = "@"
$ws.Range("U9").Value = "28.016994741766062"
$ws.Range("U9").Style = "Normal"
$ws.Range("V9").NumberFormat = "@"
$ws.Range("V9").Value = "-58.29952788916886"
$ws.Range("V9").Style = "Normal"
$ws.Range("W9").NumberFormat = "@"
$ws.Range("W9").Value = "43.47624448879986"
$ws.Range("W9").Style = "Normal"
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "-24.366290878913635"
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "-47.136113688328564"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "51.603847160216574"
$ws.Range("E10").Style = "Normal"
$ws.Range("R10").NumberFormat = "@"
$ws.Range("R10").Value = "-92.76822839730893"
$ws.Range("R10").Style = "Normal"
$ws.Range("T10").NumberFormat = "@"
$ws.Range("T10").Value = "38.38470048171806"
$ws.Range("T10").Style = "Normal"
$ws.Range("U10").NumberFormat = "@"
$ws.Range("U10").Value = "-47.136113688328564"
$ws.Range("U10").Style = "Normal"
$ws.Range("V10").NumberFormat = "@"
$ws.Range("V10").Value = "33.29486086227083"
$ws.Range("V10").Style = "Normal"
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "0.601355136578384"
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.1107878774119335"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "39.37198318836386"
$ws.Range("E11").Style = "Normal"
$ws.Range("R11").NumberFormat = "@"
$ws.Range("R11").Value = "-41.495909305675646"
$ws.Range("R11").Style = "Normal"
$ws.Range("S11").NumberFormat = "@"
$ws.Range("S11").Value = "38.12753180916628"
$ws.Range("S11").Style = "Normal"
$ws.Range("T11").NumberFormat = "@"
$ws.Range("T11").Value = "41.679754097165315"
$ws.Range("T11").Style = "Normal"
$ws.Range("U11").NumberFormat = "@"
$ws.Range("U11").Value = "-35.90595605434241"
$ws.Range("U11").Style = "Normal"
